$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their literal text representation
# (many values look numeric, e.g. "1.00", "0.999", "5.55" - without
# forcing Text format, Excel's COM layer would silently coerce them to
# numbers/doubles and we'd lose the original formatting such as
# trailing zeros or the "." thousands separators used by this sheet).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "63.199.10"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "2.615.80"
$ws.Range("E3").Value = "  +1.85%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "586.15"
$ws.Range("E5").Value = "  +2.58%  "
$ws.Range("D6").Value = "148.59"
$ws.Range("E6").Value = "  +0.91%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "0.601"
$ws.Range("E8").Value = "  +3.24%  "
$ws.Range("D9").Value = "0.109"
$ws.Range("E9").Value = "  +3.91%  "
$ws.Range("D10").Value = "5.65"
$ws.Range("E10").Value = "  +1.27%  "
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").Value = "0.357"
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("D13").Value = "27.67"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").Value = "3.082.40"
$ws.Range("E14").Value = "  +1.84%  "
$ws.Range("D15").Value = "63.046.35"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").Value = "0.0000149"
$ws.Range("E16").Value = "  +3.60%  "
$ws.Range("D17").Value = "2.696.94"
$ws.Range("E17").Value = "  +5.93%  "
$ws.Range("D18").Value = "11.45"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").Value = "346.05"
$ws.Range("E19").Value = "  +3.19%  "
$ws.Range("D20").Value = "4.44"
$ws.Range("E20").Value = "  +1.96%  "
$ws.Range("D21").Value = "6.85"
$ws.Range("E21").Value = "  -0.46%  "
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "5.55"
$ws.Range("E23").Value = "  -3.93%  "
$ws.Range("D24").Value = "66.94"
$ws.Range("E24").Value = "  +2.31%  "
$ws.Range("D25").Value = "2.734.08"
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "1.61"
$ws.Range("E27").Value = "  -0.75%  "
$ws.Range("D28").Value = "8.17"
$ws.Range("E28").Value = "  +10.12%  "
$ws.Range("D29").Value = "8.51"
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("B31").Value = "SuiNetwork"
$ws.Range("C31").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D31").Value = "1.49"
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("D32").Value = "1.97"
$ws.Range("E32").Value = "  +5.85%  "
$ws.Range("D33").Value = "0.0₃0835"
$ws.Range("E33").Value = "  +1.02%  "
$ws.Range("D34").Value = "467.82"
$ws.Range("E34").Value = "  +13.78%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "176.38"
$ws.Range("E35").Value = "  +1.03%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "1.62"
$ws.Range("E36").Value = "  +4.12%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "0.408"
$ws.Range("E38").Value = "  +1.78%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").Value = "19.37"
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D40").Value = "4.61"
$ws.Range("E40").Value = "  +5.27%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "1.73"
$ws.Range("E42").Value = "  -2.08%  "
$ws.Range("D43").Value = "156.45"
$ws.Range("E43").Value = "  +2.01%  "
$ws.Range("D44").Value = "3.83"
$ws.Range("E44").Value = "  +1.43%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "0.638"
$ws.Range("E45").Value = "  +4.90%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "20.96"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "0.0553"
$ws.Range("E47").Value = "  +3.99%  "
$ws.Range("D48").Value = "0.0977"
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("D49").Value = "0.0241"
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "18.96"
$ws.Range("E50").Value = "  +2.25%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "1.75"
$ws.Range("E51").Value = "  -1.20%  "
